$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: add I1 = "I0" and J1 = "IF", matching the style of the
# existing header cells (bold, centered, bordered -> style index 1).
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-16: new numeric values for columns I (I0) and J (IF).
# Both columns carry identical values per row.
$values = @{
    2  = 7
    3  = 8
    4  = 7
    5  = 8
    6  = 8
    7  = 5
    8  = 9
    9  = 9
    10 = 7
    11 = 9
    12 = 7
    13 = 4
    14 = 6
    15 = 5
    16 = 7
}

foreach ($row in $values.Keys) {
    $v = $values[$row]
    $ws.Cells.Item($row, 9).Value = $v
    $ws.Cells.Item($row, 10).Value = $v
}
